$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose refreshed values look like plain decimals need to stay
# text (matching the site-scraped "12.79"-style strings), so pre-format them
# as Text before writing - otherwise Excel would silently coerce them to numbers
# and drop significant trailing zeros (e.g. "0.4500" -> 0.45).
foreach ($area in $ws.Range("D5,D7,D8,D9,D10,D11,D12,D14,D15,D16,D17,D19,D20,D21,D22,D24,D25,D27,D28,D29,D30,D32,D33,D34,D36,D37,D38,D39,D40,D41,D42,D44,D45,D46,D47,D49,D50,D51").Areas) {
    $area.NumberFormat = "@"
}

$ws.Range("D2").Value = "30.606.85"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "2.111.38"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +1.13%  "
$ws.Range("D5").Value = "338.75"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D7").Value = "0.5246"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").Value = "0.4500"
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").Value = "53.27"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").Value = "0.08987"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").Value = "1.167"
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").Value = "24.34"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").Value = "2.121.17"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("D14").Value = "6.770"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").Value = "8.053"
$ws.Range("E15").Value = "  +3.02%  "
$ws.Range("D16").Value = "97.75"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").Value = "0.00001160"
$ws.Range("E17").Value = "  +2.46%  "
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").Value = "0.06698"
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("D20").Value = "19.29"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("D21").Value = "1.011"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").Value = "6.316"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "30.708.17"
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("D24").Value = "12.77"
$ws.Range("E24").Value = "  +2.95%  "
$ws.Range("D25").Value = "2.386"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("D26").Value = "2.360.05"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").Value = "22.32"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").Value = "165.17"
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("D29").Value = "2.531"
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("D30").Value = "135.30"
$ws.Range("E30").Value = "  +1.77%  "
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").Value = "0.1072"
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "1.636"
$ws.Range("E33").Value = "  -1.84%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "6.349"
$ws.Range("E34").Value = "  +2.97%  "
$ws.Range("E35").Value = "  +0.53%  "
$ws.Range("D36").Value = "10.25"
$ws.Range("E36").Value = "  -3.03%  "
$ws.Range("D37").Value = "5.873"
$ws.Range("E37").Value = "  +5.02%  "
$ws.Range("D38").Value = "0.02648"
$ws.Range("E38").Value = "  +2.51%  "
$ws.Range("D39").Value = "0.06813"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "0.2311"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("D41").Value = "12.56"
$ws.Range("E41").Value = "  -1.58%  "
$ws.Range("D42").Value = "0.6860"
$ws.Range("E43").Value = "  +0.60%  "
$ws.Range("D44").Value = "14.89"
$ws.Range("E44").Value = "  +6.06%  "
$ws.Range("D45").Value = "0.6420"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").Value = "2.306"
$ws.Range("D47").Value = "0.00000000369"
$ws.Range("E47").Value = "  +11.06%  "
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("D49").Value = "1.251"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("D50").Value = "82.68"
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("D51").Value = "0.07290"
$ws.Range("E51").Value = "  +2.76%  "
